$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (including the date number format) from G7 down to G8
$ws.Range("G7").Copy($ws.Range("G8"))

$ws.Range("A8").Value = 9815.06
$ws.Range("B8").Value = 9697.7199999999993
$ws.Range("C8").Value = 309.02999999999997
$ws.Range("D8").Value = 305.29000000000002
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = -1.21
$ws.Range("G8").Value = 42608.617939814816
$ws.Range("H8").Value = $true
